$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.525.20'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.25%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.604.47'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.02%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.08'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.06%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.527'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +8.60%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.997'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.19%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '26.90'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +9.64%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '43.43'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.21%  '
$ws.Range("E10").Value = '  +2.18%  '
$ws.Range("E11").Value = '  +2.63%  '
$ws.Range("E12").Value = '  +1.85%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.834.36'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.07%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.592.14'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.20%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '29.553.95'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.34%  '
$ws.Range("E16").Value = '  +4.42%  '
$ws.Range("E17").Value = '  +2.65%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '63.49'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.51%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '243.00'
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.62'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.38%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0692'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.06%  '
$ws.Range("E22").Value = '  -0.11%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.25%  '
$ws.Range("E24").Value = '  +2.41%  '
$ws.Range("E25").Value = '  -0.07%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '154.58'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.22%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.110'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +6.57%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.32'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.92%  '
$ws.Range("E29").Value = '  +2.73%  '
$ws.Range("E30").Value = '  -0.15%  '
$ws.Range("E31").Value = '  +2.76%  '
$ws.Range("E32").Value = '  +0.04%  '
$ws.Range("E33").Value = '  +1.89%  '
$ws.Range("E34").Value = '  +3.77%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.416.15'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.62%  '
$ws.Range("E36").Value = '  -2.27%  '
$ws.Range("E37").Value = '  +2.36%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.81'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.90%  '
$ws.Range("E39").Value = '  +1.41%  '
$ws.Range("E40").Value = '  +1.85%  '
$ws.Range("E41").Value = '  +3.33%  '
$ws.Range("E42").Value = '  +1.02%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.997'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.21%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '52.76'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +21.44%  '
$ws.Range("E45").Value = '  +2.15%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0474'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.36%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '65.75'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.72%  '
$ws.Range("E48").Value = '  -0.36%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.745.36'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.63%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '86.34'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.43%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.832'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.31%  '
